$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 569
$ws.Range("F5").Value = 6341
$ws.Range("F6").Value = 714
$ws.Range("F8").Value = 65
$ws.Range("F9").Value = 624
$ws.Range("F10").Value = 313
$ws.Range("F11").Value = 190
$ws.Range("F12").Value = 676
$ws.Range("F13").Value = 6
$ws.Range("F14").Value = 1148
$ws.Range("F16").Value = 403
$ws.Range("F17").Value = 47
$ws.Range("F18").Value = 18
$ws.Range("F19").Value = 1413
$ws.Range("F20").Value = 656
$ws.Range("F21").Value = 373
$ws.Range("F22").Value = 393
$ws.Range("F25").Value = 124
$ws.Range("F26").Value = 2185
$ws.Range("F27").Value = 240
$ws.Range("F28").Value = 86
$ws.Range("F29").Value = 389
$ws.Range("F30").Value = 57
$ws.Range("F31").Value = 3525

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 164
$ws.Range("F8").Value = 701
$ws.Range("F14").Value = 636
$ws.Range("F20").Value = 4082
$ws.Range("F24").Value = 184
$ws.Range("F28").Value = 207
$ws.Range("F32").Value = 1591

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 1183
$ws.Range("F7").Value = 1567
$ws.Range("F9").Value = 130
$ws.Range("F11").Value = 743

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1183
$ws.Range("F5").Value = 1567
$ws.Range("F7").Value = 130
$ws.Range("F8").Value = 743
$ws.Range("F9").Value = 569
$ws.Range("F12").Value = 6341
$ws.Range("F15").Value = 714
$ws.Range("F17").Value = 701
$ws.Range("F18").Value = 624
$ws.Range("F19").Value = 190
$ws.Range("F20").Value = 676
$ws.Range("F25").Value = 1148
$ws.Range("F26").Value = 403
$ws.Range("F29").Value = 47
$ws.Range("F30").Value = 18
$ws.Range("F31").Value = 1413
$ws.Range("F34").Value = 656
$ws.Range("F35").Value = 373
$ws.Range("F36").Value = 393
$ws.Range("F39").Value = 184
$ws.Range("F43").Value = 207
$ws.Range("F45").Value = 1591
$ws.Range("F46").Value = 240
$ws.Range("F47").Value = 86
$ws.Range("F48").Value = 389
$ws.Range("F49").Value = 57
$ws.Range("F50").Value = 3525
